# Auto-generated edit script applying market-price / profit recalculation updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 154.75
$ws.Range("I2").Value = 154.75
$ws.Range("K2").Value = 154.75
$ws.Range("M2").Value = -41.75
$ws.Range("H17").Value = 449.25
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = ""
$ws.Range("H44").Value = 18691.666
$ws.Range("I44").Value = 2750
$ws.Range("J44").Value = 26662.5
$ws.Range("K44").Value = 2750
$ws.Range("L44").Value = 26662.5
$ws.Range("M44").Value = -2288
$ws.Range("N44").Value = -27586.5
$ws.Range("H82").Value = 21833.334
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 21833.334
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 65500.00199999999
$ws.Range("M82").Value = ""
$ws.Range("N82").Value = -66312.00199999999
$ws.Range("H85").Value = 21833.334
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 21833.334
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 65500.00199999999
$ws.Range("M85").Value = ""
$ws.Range("N85").Value = -68308.00199999999
$ws.Range("H137").Value = 2749.3462
$ws.Range("I137").Value = 1499.0555
$ws.Range("J137").Value = 5562.5
$ws.Range("K137").Value = 4497.166499999999
$ws.Range("L137").Value = 16687.5
$ws.Range("M137").Value = -1947.166499999999
$ws.Range("N137").Value = -21787.5
$ws.Range("H138").Value = 2348.0762
$ws.Range("I138").Value = 2245.963
$ws.Range("J138").Value = 2390.4922
$ws.Range("K138").Value = 6737.889000000001
$ws.Range("L138").Value = 7171.4766
$ws.Range("M138").Value = -1597.889000000001
$ws.Range("N138").Value = -17451.4766
$ws.Range("H141").Value = 2711.923
$ws.Range("I141").Value = 4791.25
$ws.Range("J141").Value = 2474.2856
$ws.Range("K141").Value = 14373.75
$ws.Range("L141").Value = 7422.8568
$ws.Range("M141").Value = -9193.75
$ws.Range("N141").Value = -17782.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = ""
$ws.Range("H32").Value = 45951.46
$ws.Range("I32").Value = 12897.345
$ws.Range("K32").Value = 12897.345
$ws.Range("M32").Value = -12610.345
$ws.Range("H107").Value = 11600
$ws.Range("J107").Value = 11600
$ws.Range("L107").Value = 11600
$ws.Range("N107").Value = -19280
$ws.Range("H133").Value = 31497.5
$ws.Range("J133").Value = 31497.5
$ws.Range("L133").Value = 31497.5
$ws.Range("N133").Value = -36557.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = ""
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H81").Value = 19260
$ws.Range("J81").Value = 19260
$ws.Range("L81").Value = 19260
$ws.Range("N81").Value = -21382
$ws.Range("H84").Value = 19260
$ws.Range("J84").Value = 19260
$ws.Range("L84").Value = 57780
$ws.Range("N84").Value = -68388
$ws.Range("H99").Value = 1690.5454
$ws.Range("I99").Value = 1055.5
$ws.Range("K99").Value = 1055.5
$ws.Range("M99").Value = 442.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 8162
$ws.Range("J14").Value = 8000
$ws.Range("L14").Value = 8000
$ws.Range("N14").Value = -8340
$ws.Range("H50").Value = 9207.200000000001
$ws.Range("J50").Value = 9207.200000000001
$ws.Range("L50").Value = 9207.200000000001
$ws.Range("N50").Value = -10457.2
$ws.Range("H60").Value = 24813.715
$ws.Range("J60").Value = 24813.715
$ws.Range("L60").Value = 24813.715
$ws.Range("N60").Value = -25835.715
$ws.Range("H68").Value = 17600.666
$ws.Range("J68").Value = 17600.666
$ws.Range("L68").Value = 17600.666
$ws.Range("N68").Value = -19098.666
$ws.Range("H71").Value = 17600.666
$ws.Range("J71").Value = 17600.666
$ws.Range("L71").Value = 52801.99800000001
$ws.Range("N71").Value = -60289.99800000001
$ws.Range("H88").Value = 31594.375
$ws.Range("J88").Value = 31594.375
$ws.Range("L88").Value = 31594.375
$ws.Range("N88").Value = -32406.375
$ws.Range("H91").Value = 31594.375
$ws.Range("J91").Value = 31594.375
$ws.Range("L91").Value = 31594.375
$ws.Range("N91").Value = -34402.375
$ws.Range("H107").Value = 898.1111
$ws.Range("I107").Value = 936.38464
$ws.Range("K107").Value = 936.38464
$ws.Range("M107").Value = 983.61536
$ws.Range("H109").Value = 11450
$ws.Range("J109").Value = 11450
$ws.Range("L109").Value = 11450
$ws.Range("N109").Value = -13530
$ws.Range("H138").Value = 39730
$ws.Range("J138").Value = 39730
$ws.Range("L138").Value = 39730
$ws.Range("N138").Value = -50010
$ws.Range("H139").Value = 54499
$ws.Range("J139").Value = 54499
$ws.Range("L139").Value = 54499
$ws.Range("N139").Value = -64779

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3069.7778
$ws.Range("J58").Value = 3069.7778
$ws.Range("L58").Value = 9209.3334
$ws.Range("N58").Value = -9465.3334
$ws.Range("H80").Value = 5528.8823
$ws.Range("I80").Value = 3820.2
$ws.Range("J80").Value = 6240.8335
$ws.Range("K80").Value = 11460.6
$ws.Range("L80").Value = 18722.5005
$ws.Range("M80").Value = -10524.6
$ws.Range("N80").Value = -20594.5005
$ws.Range("H83").Value = 5528.8823
$ws.Range("I83").Value = 3820.2
$ws.Range("J83").Value = 6240.8335
$ws.Range("K83").Value = 34381.8
$ws.Range("L83").Value = 56167.5015
$ws.Range("M83").Value = -29701.8
$ws.Range("N83").Value = -65527.5015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 13487.75
$ws.Range("I126").Value = 2998.75
$ws.Range("K126").Value = 8996.25
$ws.Range("M126").Value = -6526.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = ""
$ws.Range("H133").Value = 56500
$ws.Range("J133").Value = 56500
$ws.Range("L133").Value = 56500
$ws.Range("N133").Value = -61560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10267.625
$ws.Range("I32").Value = 750
$ws.Range("J32").Value = 13440.167
$ws.Range("K32").Value = 750
$ws.Range("L32").Value = 13440.167
$ws.Range("M32").Value = -433
$ws.Range("N32").Value = -14074.167
$ws.Range("H109").Value = 25000
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27774
$ws.Range("H138").Value = 35249
$ws.Range("J138").Value = 35249
$ws.Range("L138").Value = 35249
$ws.Range("N138").Value = -45529
